$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM data: the dataset was recomputed and now only contains the
# FAPs -> Cxcl5/Cxcr1 -> MuSCs interaction (the "ECs" sending/target
# cluster rows are gone entirely). Update the surviving row's values to
# the freshly-derived specificities (recalculated against the smaller
# dataset, so several specificity columns become 1), then drop the
# rows that no longer exist.

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl5"
$ws.Range("C2").Value = "Cxcr1"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05862133333333333
$ws.Range("H2").Value = 0.175864
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.000484
$ws.Range("N2").Value = 0.001452
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.00002837272533333333
$ws.Range("R2").Value = 0.000255354528
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Drop the rows that described the now-removed "ECs" cluster
# interactions (rows 3-5 in the original sheet).
$ws.Rows("3:5").Delete()
